$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")

$ws.Range("D2").Value = "26.647.37"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").Value = "1.686.44"
$ws.Range("E3").Value = "  +3.39%  "
$helper.Formula = "=""1.001"""
$helper.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  -0.07%  "
$helper.Formula = "=""217.27"""
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +3.74%  "
$helper.Formula = "=""0.5338"""
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +4.51%  "
$helper.Formula = "=""0.06432"""
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +3.43%  "
$helper.Formula = "=""21.65"""
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +7.23%  "
$helper.Formula = "=""0.07804"""
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("D12").Value = "1.696.21"
$ws.Range("E12").Value = "  +3.68%  "
$helper.Formula = "=""4.498"""
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +3.55%  "
$helper.Formula = "=""0.5612"""
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("D15").Value = "0.0₅8447"
$ws.Range("E15").Value = "  +6.68%  "
$helper.Formula = "=""66.09"""
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("D17").Value = "26.678.02"
$ws.Range("E17").Value = "  +2.95%  "
$helper.Formula = "=""1.001"""
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -0.10%  "
$helper.Formula = "=""4.796"""
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +4.13%  "
$helper.Formula = "=""195.24"""
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +5.98%  "
$helper.Formula = "=""10.40"""
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +3.90%  "
$helper.Formula = "=""6.375"""
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +5.28%  "
$helper.Formula = "=""1.002"""
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -0.08%  "
$helper.Formula = "=""144.48"""
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.77%  "
$helper.Formula = "=""0.1284"""
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +6.23%  "
$helper.Formula = "=""7.477"""
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +2.03%  "
$helper.Formula = "=""16.30"""
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +5.32%  "
$helper.Formula = "=""1.435"""
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +4.92%  "
$helper.Formula = "=""0.06155"""
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +4.24%  "
$helper.Formula = "=""1.278"""
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +3.01%  "
$helper.Formula = "=""3.613"""
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +7.91%  "
$helper.Formula = "=""3.467"""
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +3.70%  "
$helper.Formula = "=""1.700"""
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +6.18%  "
$ws.Range("E34").Value = "  +4.60%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$helper.Formula = "=""2.796"""
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$helper.Formula = "=""2.423"""
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +1.66%  "
$helper.Formula = "=""0.5739"""
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -0.34%  "
$helper.Formula = "=""0.01648"""
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +3.44%  "
$helper.Formula = "=""6.033"""
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +7.66%  "
$ws.Range("D40").Value = "1.070.41"
$ws.Range("E40").Value = "  +5.96%  "
$helper.Formula = "=""0.8638"""
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +3.65%  "
$helper.Formula = "=""1.001"""
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +0.02%  "
$helper.Formula = "=""100.40"""
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "1.837.10"
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("D45").Value = "0.0₈110"
$ws.Range("E45").Value = "  +0.17%  "
$helper.Formula = "=""57.28"""
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +5.72%  "
$helper.Formula = "=""8.170"""
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("E48").Value = "  +0.59%  "
$helper.Formula = "=""0.05219"""
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +0.91%  "
$helper.Formula = "=""6.079"""
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +5.36%  "
$helper.Formula = "=""0.4240"""
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.43%  "

$helper.Value = ""
$excel.CutCopyMode = $false
